$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: relabel rows, drop the old header text and shift the "max+N" labels up ---
$ws.Range("A1").Value = $null

$ws.Range("A2").Value = "max+1"
$ws.Range("A3").Value = "max+2"
$ws.Range("A4").Value = "max+3"
$ws.Range("A5").Value = "max+4"
$ws.Range("A6").Value = "max+5"
$ws.Range("A7").Value = "max+6"
$ws.Range("A8").Value = "max+7"
$ws.Range("A9").Value = "max+8"
$ws.Range("A10").Value = "max+9"
$ws.Range("A11").Value = "max+10"

# Row 12 is no longer part of the table.
$ws.Range("A12:E12").ClearContents()

# --- Column C: retune values, only rows 2-3 keep data now ---
$ws.Range("C2").Value = 8
$ws.Range("C3").Value = 0
$ws.Range("C4:C10").ClearContents()

# --- Column D: retune values, rows 2-6 keep data, row 7 cleared ---
$ws.Range("D2").Value = 31
$ws.Range("D3").Value = 10
$ws.Range("D4").Value = 2
$ws.Range("D5").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("D7").ClearContents()

# --- Column E: retune values for rows 2-11 ---
$ws.Range("E2").Value = 85
$ws.Range("E3").Value = 58
$ws.Range("E4").Value = 33
$ws.Range("E5").Value = 18
$ws.Range("E6").Value = 8
$ws.Range("E7").Value = 3
$ws.Range("E8").Value = 2
$ws.Range("E9").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("E11").Value = 0

# --- Sheet view: selection moves from I5 to R6 ---
$ws.Range("R6").Select() | Out-Null

# --- Book view (window geometry) ---
$win = $wb.Windows.Item(1)
$win.Left = -103
$win.Top = -103
$win.Width = 21806
$win.Height = 13886
